# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (column H) and
# "Correspond Handback DateTime" (column K) values on the zh-cn and
# de-de sheets for the 34231130-...md file (row 2), reflecting a fresh
# handback report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-20 12:50:56"
$wsZhCn.Range("K2").Value = "2016-08-20 12:51:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-20 12:51:00"
$wsDeDe.Range("K2").Value = "2016-08-20 12:51:20"
